$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: mark the "Tavolo di gioco / Realizzazione interfaccia basilare" task as
# finished on day 1 (remaining effort "/")
$ws.Range("E5").Value = "/"

# Row 6: fill in the daily remaining-effort burndown for the
# "Interazione di pesca e scarto da parte dell'utente" task (started at 4,
# counted down 4,3,2,1 then finished "/")
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = "/"
